$wb = $excel.ActiveWorkbook

# Sheets: 1 = CONTEÚDO, 2 = EXEMPLO, 3 = EXERCICIO
$wsConteudo = $wb.Worksheets.Item(1)
$wsExemplo  = $wb.Worksheets.Item(2)

# --- Values for the right-hand "inputs" column ---
$wsExemplo.Range("I10").Value = 50
$wsExemplo.Range("I11").Value = 100
$wsExemplo.Range("I12").Value = 100
$wsExemplo.Range("I14").Value = 150

# Writing to I11 (a cell with no explicit row/cell style) makes the engine
# recompute an autofit row height; put row 11 back to its original height.
$wsExemplo.Rows(11).RowHeight = 3.95

# --- Formulas for the left-hand "results" column ---
$wsExemplo.Range("E10").Formula = "=I12*I10*I14"
$wsExemplo.Range("E12").Formula = "=(I14/I10)/I12"
$wsExemplo.Range("E14").Formula = "=I10+I12+I14"
$wsExemplo.Range("E16").Formula = "=SQRT(I12)"
$wsExemplo.Range("E18").Formula = "=PI()"
$wsExemplo.Range("E20").Formula = "=I10^I12"

# E12 (division result) gets a 2-decimal number format (new style in cellXfs)
$wsExemplo.Range("E12").NumberFormat = "0.00"

# --- Column E width on EXEMPLO ---
$wsExemplo.Columns("E").ColumnWidth = 11.1

# --- Selection / active sheet changes ---
$wsConteudo.Range("F11").Select()
$wsExemplo.Activate()
$wsExemplo.Range("G26").Select()
